# Added iMethodInterceptor to disable test case at run time
# -> add a new "RUNMANAGER" worksheet (after "testing") that drives which
#    tests execute, at what priority, and how many times.

$wb = $excel.ActiveWorkbook

$testingSheet = $wb.Worksheets.Item(1)

# New sheet goes right after "testing".
$runManager = $wb.Worksheets.Add($null, $testingSheet)
$runManager.Name = "RUNMANAGER"

# Columns D (priority) and E (count) are stored as text.
$runManager.Range("D1:E3").NumberFormat = "@"

# Header row.
$runManager.Range("A1").Value = "testName"
$runManager.Range("B1").Value = "testDescripton"
$runManager.Range("C1").Value = "execute"
$runManager.Range("D1").Value = "priority"
$runManager.Range("E1").Value = "count"

# Row 2.
$runManager.Range("A2").Value = "loginLogoutTest"
$runManager.Range("B2").Value = "Test the login and logout of the app"
$runManager.Range("C3").Value = "yes"

# Row 3.
$runManager.Range("A3").Value = "newTest"
$runManager.Range("B3").Value = "This is second test"
$runManager.Range("C2").Value = "no"

$runManager.Range("D2").Value = "'1"
$runManager.Range("D3").Value = "'2"
$runManager.Range("E2").Value = "'1"
$runManager.Range("E3").Value = "'1"

# Match the author's manual column resize on the new sheet.
[void]$runManager.Columns("A:B").AutoFit()

$runManager.PageSetup.Orientation = 1
$runManager.PageSetup.PaperSize = 9

# Restore the selection on "testing" (no longer the active tab) ...
$testingSheet.Activate()
[void]$testingSheet.Range("H11").Select()

# ... and leave RUNMANAGER as the active sheet/selection.
$runManager.Activate()
[void]$runManager.Range("E10").Select()
